$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.933.58"
$ws.Range("D3").Value = "2.225.17"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -1.80%  "
$ws.Range("D5").Value = "'299.24"
$ws.Range("E5").Value = "  -3.11%  "
$ws.Range("D6").Value = "'90.33"
$ws.Range("E6").Value = "  -4.54%  "
$ws.Range("D7").Value = "'0.554"
$ws.Range("E7").Value = "  -3.49%  "
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").Value = "'0.490"
$ws.Range("E9").Value = "  -6.72%  "
$ws.Range("D10").Value = "'32.95"
$ws.Range("E10").Value = "  -5.53%  "
$ws.Range("D11").Value = "'0.0778"
$ws.Range("E11").Value = "  -3.80%  "
$ws.Range("D12").Value = "'6.93"
$ws.Range("E12").Value = "  -4.47%  "
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "2.565.22"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "2.226.73"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").Value = "'13.44"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "'0.774"
$ws.Range("E17").Value = "  -7.54%  "
$ws.Range("D18").Value = "43.803.17"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").Value = "0.0₃0902"
$ws.Range("E19").Value = "  -5.96%  "
$ws.Range("E20").Value = "  -6.83%  "
$ws.Range("D21").Value = "'11.25"
$ws.Range("E21").Value = "  -6.97%  "
$ws.Range("D22").Value = "'64.64"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "'236.40"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").Value = "'2.80"
$ws.Range("E24").Value = "  -6.10%  "
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "'1.86"
$ws.Range("E26").Value = "  -6.23%  "
$ws.Range("D27").Value = "'38.16"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").Value = "'9.31"
$ws.Range("E29").Value = "  -5.23%  "
$ws.Range("D30").Value = "'19.23"
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("D31").Value = "'150.89"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").Value = "'5.40"
$ws.Range("E32").Value = "  -9.76%  "
$ws.Range("D33").Value = "'0.0751"
$ws.Range("E33").Value = "  -6.40%  "
$ws.Range("D34").Value = "'2.51"
$ws.Range("E34").Value = "  -5.30%  "
$ws.Range("E35").Value = "  -3.81%  "
$ws.Range("D36").Value = "'2.82"
$ws.Range("E36").Value = "  -10.34%  "
$ws.Range("E37").Value = "  -6.59%  "
$ws.Range("E38").Value = "  -5.92%  "
$ws.Range("D39").Value = "'0.0300"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "'3.20"
$ws.Range("E40").Value = "  -6.74%  "
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("D42").Value = "'13.15"
$ws.Range("E42").Value = "  -8.30%  "
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").Value = "1.836.45"
$ws.Range("E44").Value = "  +4.81%  "
$ws.Range("D45").Value = "'1.78"
$ws.Range("E45").Value = "  +12.61%  "
$ws.Range("D46").Value = "'0.181"
$ws.Range("E46").Value = "  -6.27%  "
$ws.Range("D47").Value = "'67.23"
$ws.Range("E47").Value = "  -4.39%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'73.26"
$ws.Range("E48").Value = "  -9.03%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'93.96"
$ws.Range("E49").Value = "  -5.53%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'13.99"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").Value = "2.446.96"
$ws.Range("E51").Value = "  -1.06%  "
